$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07794266666666667
$ws.Range("H2").Value = 0.233828
$ws.Range("I2").Value = 0.002827880818927331
$ws.Range("J2").Value = 0.00282788081892733
$ws.Range("Q2").Value = 0.897683722736889
$ws.Range("R2").Value = 8.079153504632
$ws.Range("S2").Value = 0.002827880818927331
$ws.Range("T2").Value = 0.00282788081892733

# Row 3
$ws.Range("I3").Value = 0.9151728997907317
$ws.Range("J3").Value = 0.9151728997907316
$ws.Range("S3").Value = 0.9151728997907317
$ws.Range("T3").Value = 0.9151728997907316

# Row 4
$ws.Range("G4").Value = 2.260080333333333
$ws.Range("H4").Value = 6.780241
$ws.Range("I4").Value = 0.08199921939034102
$ws.Range("J4").Value = 0.08199921939034102
$ws.Range("Q4").Value = 26.02986803091711
$ws.Range("R4").Value = 234.268812278254
$ws.Range("S4").Value = 0.08199921939034102
$ws.Range("T4").Value = 0.08199921939034102
